# Refactor data parsing logic
# Appends a new row (row 19) of parsed data to each of the four worksheets,
# mirroring the structure of the existing rows (A: timestamp, B-E: hex byte
# strings, F-I: decimal decodings).

$wb = $excel.ActiveWorkbook

$rows = @{
    "ROW35-FE-LIFTER"  = @{
        A = [double]"45733.7684509375"
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x86"
        E = "0xd"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 390
        I = 13
    }
    "ROW35-MID-LIFTER" = @{
        A = [double]"45733.61935922454"
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x86"
        E = "0xe"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 390
        I = 14
    }
    "ROW02-FE-LIFTER"  = @{
        A = [double]"45733.76750971065"
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x86"
        E = "0x3"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 390
        I = 3
    }
    "ROW02-MID-LIFTER" = @{
        A = [double]"45733.82578517361"
        B = "0x01,0x90"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x01,0x86"
        E = "0x3"
        F = 400
        G = [double]"9.85046333984776e+23"
        H = 390
        I = 3
    }
}

foreach ($ws in $wb.Worksheets) {
    $data = $rows[$ws.Name]
    if ($null -eq $data) {
        continue
    }

    $newRow = $ws.UsedRange.Rows.Count + 1

    $ws.Cells.Item($newRow, 1).Value = $data.A
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat
    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = $data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
